# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (B16:G27) alternates rows between the
# two workers (YORLENE MONTIEL ALVAREZ / GLADYS ESTHER CERVANTES VARGAS),
# one row per worker per mora period (1901-1906). The data is being
# reordered so the table is sorted by ascending period (1901..1906),
# interleaving both workers for each period, instead of the previous
# ordering (descending period, grouped by worker).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1901
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1044911851"
$ws.Range("D16").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

# Row 17: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1901
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45497273"
$ws.Range("D17").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1200000

# Row 18: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1902
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1044911851"
$ws.Range("D18").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 828116

# Row 19: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1902
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45497273"
$ws.Range("D19").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 48000
$ws.Range("G19").Value = 1200000

# Row 20: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1902 (unchanged)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1044911851"
$ws.Range("D20").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E20").Value = "1902"
$ws.Range("F20").Value = 33125
$ws.Range("G20").Value = 828116

# Row 21: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1903
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45497273"
$ws.Range("D21").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E21").Value = "1903"
$ws.Range("F21").Value = 48000
$ws.Range("G21").Value = 1200000

# Row 22: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1904
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1044911851"
$ws.Range("D22").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E22").Value = "1904"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116

# Row 23: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1904
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "45497273"
$ws.Range("D23").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E23").Value = "1904"
$ws.Range("F23").Value = 48000
$ws.Range("G23").Value = 1200000

# Row 24: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1905
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1044911851"
$ws.Range("D24").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E24").Value = "1905"
$ws.Range("F24").Value = 33125
$ws.Range("G24").Value = 828116

# Row 25: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1905
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "45497273"
$ws.Range("D25").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E25").Value = "1905"
$ws.Range("F25").Value = 48000
$ws.Range("G25").Value = 1200000

# Row 26: CC / 1044911851 / YORLENE MONTIEL ALVAREZ / periodo 1906
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1044911851"
$ws.Range("D26").Value = "YORLENE MONTIEL ALVAREZ"
$ws.Range("E26").Value = "1906"
$ws.Range("F26").Value = 18771
$ws.Range("G26").Value = 828116

# Row 27: CC / 45497273 / GLADYS ESTHER CERVANTES VARGAS / periodo 1906
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "45497273"
$ws.Range("D27").Value = "GLADYS ESTHER CERVANTES VARGAS"
$ws.Range("E27").Value = "1906"
$ws.Range("F27").Value = 27200
$ws.Range("G27").Value = 1200000
